$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "301.53"
Set-TextValue "E2" "2.59%"
Set-TextValue "D3" "32.02"
Set-TextValue "E3" "3.03%"
Set-TextValue "D4" "5.021"
Set-TextValue "E4" "1.71%"
Set-TextValue "E5" "6.10%"
Set-TextValue "D6" "2.313"
Set-TextValue "E6" "0.99%"
Set-TextValue "D7" "7.976"
Set-TextValue "E7" "3.92%"
Set-TextValue "D8" "0.9319"
Set-TextValue "E8" "2.13%"
Set-TextValue "D9" "0.1021"
Set-TextValue "E9" "25.11%"
Set-TextValue "D10" "0.1772"
Set-TextValue "E10" "5.20%"
Set-TextValue "D11" "0.08468"
Set-TextValue "E11" "3.04%"
Set-TextValue "D12" "0.03331"
Set-TextValue "E12" "7.30%"
Set-TextValue "D13" "0.09892"
Set-TextValue "E13" "-1.65%"
Set-TextValue "D14" "0.001472"
Set-TextValue "E14" "-2.98%"
Set-TextValue "D15" "0.005721"
Set-TextValue "E15" "-0.35%"
Set-TextValue "D17" "3.861"
Set-TextValue "E17" "2.10%"
Set-TextValue "D18" "2.192"
Set-TextValue "E18" "5.34%"
Set-TextValue "D19" "0.3362"
Set-TextValue "E19" "1.00%"
Set-TextValue "D20" "0.1340"
Set-TextValue "E20" "2.80%"
Set-TextValue "D21" "4.307"
Set-TextValue "E21" "8.67%"
Set-TextValue "E22" "-0.96%"
Set-TextValue "D23" "0.04633"
Set-TextValue "E23" "1.99%"
Set-TextValue "E24" "0.62%"
Set-TextValue "D25" "0.004384"
Set-TextValue "E25" "0.94%"
Set-TextValue "D26" "0.0001293"
Set-TextValue "E26" "-0.54%"
Set-TextValue "D27" "0.0003369"
Set-TextValue "E27" "-0.80%"
Set-TextValue "D39" "0.01705"
Set-TextValue "E39" "6.57%"
Set-TextValue "D40" "0.04764"
Set-TextValue "E40" "7.44%"
Set-TextValue "D41" "0.007711"
Set-TextValue "E41" "4.86%"
Set-TextValue "D42" "0.009760"
Set-TextValue "E42" "11.55%"
Set-TextValue "D43" "0.1404"
Set-TextValue "E43" "5.81%"
Set-TextValue "D44" "0.002073"
Set-TextValue "E44" "-1.73%"
Set-TextValue "D45" "0.009668"
Set-TextValue "E45" "5.08%"
Set-TextValue "D46" "0.00006093"
Set-TextValue "E46" "2.09%"
Set-TextValue "D47" "0.00000000745"
Set-TextValue "E47" "-0.79%"
Set-TextValue "D48" "2.794"
Set-TextValue "E48" "24.67%"
Set-TextValue "D49" "0.001985"
Set-TextValue "E49" "-31.51%"
Set-TextValue "D50" "0.00002085"
Set-TextValue "E50" "-0.79%"
Set-TextValue "D51" "0.0001985"
Set-TextValue "E51" "-0.79%"
